# Fruta / hortaliza, semanal
#
# The underlying data rows (2-41) get reshuffled: each destination row ends
# up holding the values that used to live in a different source row. Row 1
# (headers) and rows 14/15 (unchanged) are left alone. We snapshot every
# source row's values first (so later writes don't clobber data we still
# need to read), then write the permuted rows back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol  = 20  # T

# destinationRow -> sourceRow (source row's original contents move into
# the destination row)
$map = @{
    2  = 17
    3  = 37
    4  = 23
    5  = 25
    6  = 38
    7  = 10
    8  = 13
    9  = 11
    10 = 41
    11 = 33
    12 = 39
    13 = 7
    14 = 14
    15 = 15
    16 = 4
    17 = 30
    18 = 9
    19 = 29
    20 = 2
    21 = 32
    22 = 35
    23 = 8
    24 = 19
    25 = 20
    26 = 18
    27 = 31
    28 = 12
    29 = 16
    30 = 6
    31 = 22
    32 = 34
    33 = 3
    34 = 28
    35 = 24
    36 = 27
    37 = 40
    38 = 36
    39 = 21
    40 = 5
    41 = 26
}

# Snapshot the original value of every cell in rows 2-41 before any writes.
$snapshot = @{}
for ($r = 2; $r -le 41; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# Write each destination row using the snapshot of its mapped source row.
foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $snapshot["$srcRow,$c"]
    }
}
